# Add nested route structure
# Reposition a handful of shapes on slide 1 and slide 2 (EMU -> points, 1 pt = 12700 EMU).

$p = $ppt.ActivePresentation

# NOTE: Shape.Left/Top/Width/Height are Single (float32) points under the
# hood, and EMU = round(points * 12700) truncates toward zero. Nudging by
# half an EMU before the points conversion keeps the round-trip exact.

# --- Slide 1 ---
$s1 = $p.Slides.Item(1)

# "직사각형 4" - off 3383006,2321004 / ext 5425987,2215991 -> off 3200399,2619999 / ext 5791200,2215991
$shp = $s1.Shapes.Item(3)
$shp.Left = (3200399 + 0.5) / 12700
$shp.Top = (2619999 + 0.5) / 12700
$shp.Width = (5791200 + 0.5) / 12700
$shp.Height = (2215991 + 0.5) / 12700

# --- Slide 2 ---
$s2 = $p.Slides.Item(2)

# "TextBox 7" - off 761998,2780844 -> 761998,2850293
$shp = $s2.Shapes.Item(5)
$shp.Left = (761998 + 0.5) / 12700
$shp.Top = (2850293 + 0.5) / 12700

# "사각형: 둥근 모서리 8" - off 3912243,3939712 -> 3912243,4009161
$shp = $s2.Shapes.Item(6)
$shp.Left = (3912243 + 0.5) / 12700
$shp.Top = (4009161 + 0.5) / 12700

# "TextBox 9" - off 1444903,3957680 -> 1444903,4027129
$shp = $s2.Shapes.Item(7)
$shp.Left = (1444903 + 0.5) / 12700
$shp.Top = (4027129 + 0.5) / 12700

# "TextBox 10" - off 1444902,4911753 -> 1444902,4981202
$shp = $s2.Shapes.Item(8)
$shp.Left = (1444902 + 0.5) / 12700
$shp.Top = (4981202 + 0.5) / 12700

# "사각형: 둥근 모서리 11" - off 3912243,4893784 -> 3912243,4963233
$shp = $s2.Shapes.Item(9)
$shp.Left = (3912243 + 0.5) / 12700
$shp.Top = (4963233 + 0.5) / 12700
